# Trapper repeared, arch 30/12/25
# Applies the data refresh to the YData sheet of the YahooDataOutput workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ticker rename: LAZR -> LAZRQ (row 62)
$ws.Range("A62").Value = "LAZRQ"

# Remove the stale "--" placeholder cells (Earnings Date column, rows 38 and 50)
$ws.Range("D38").ClearContents()
$ws.Range("D50").ClearContents()

# Refreshed "1-Year Target Price" values (column B)
$ws.Range("B3").Value = 22.75
$ws.Range("B8").Value = 456.802
$ws.Range("B12").Value = 118.61905
$ws.Range("B14").Value = 94.59999999999999
$ws.Range("B25").Value = 475.4725
$ws.Range("B26").Value = 73.81818
$ws.Range("B31").Value = 90.69091
$ws.Range("B32").Value = 109.495
$ws.Range("B39").Value = 51.2
$ws.Range("B48").Value = 456.802
$ws.Range("B62").Value = 1
$ws.Range("B64").Value = 657.92224

# Refreshed "Dividend Yield (%)" values (column E)
$ws.Range("E3").Value = 0.0116
$ws.Range("E8").Value = 0.0074
$ws.Range("E9").Value = 0.016900001
$ws.Range("E13").Value = 0.0211
$ws.Range("E14").Value = 0.0332
$ws.Range("E16").Value = 0.013099999
$ws.Range("E20").Value = 0.0343
$ws.Range("E23").Value = 0.0234
$ws.Range("E24").Value = 0.052800003
$ws.Range("E25").Value = 0.0226
$ws.Range("E26").Value = 0.0223
$ws.Range("E31").Value = 0.028199999
$ws.Range("E32").Value = 0.01
$ws.Range("E35").Value = 0.0107
$ws.Range("E40").Value = 0.024400001
$ws.Range("E42").Value = 0.0226
$ws.Range("E45").Value = 0.0112
$ws.Range("E48").Value = 0.0074
$ws.Range("E49").Value = 0.025
$ws.Range("E50").Value = 0.0468
$ws.Range("E51").Value = 0.0388
$ws.Range("E53").Value = 0.0342
$ws.Range("E56").Value = 0.06510000000000001
$ws.Range("E57").Value = 0.026099999
$ws.Range("E58").Value = 0.0452
$ws.Range("E59").Value = 0.025
$ws.Range("E61").Value = 0.0179
$ws.Range("E63").Value = 0.0175
$ws.Range("E66").Value = 0.0287
$ws.Range("E67").Value = 0.0531
$ws.Range("E69").Value = 0.0688
$ws.Range("E71").Value = 0.01
